$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 81; existing rows 81-117 shift down to 82-118.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new weekly record.
$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44795
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = 100112035
$ws.Cells.Item(81, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 160
$ws.Cells.Item(81, 11).Value = 23000
$ws.Cells.Item(81, 12).Value = 25000
$ws.Cells.Item(81, 13).Value = 24000
$ws.Cells.Item(81, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 2400
$ws.Cells.Item(81, 17).Value = 10
$ws.Cells.Item(81, 18).Value = "Hortaliza"
